$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Houses of Worship" policy indicator (column H) is cleared to 0
# for rows 35 through 172 (previously flagged as active with value 1).
$ws.Range("H35:H172").Value = 0
